# Insert a new data row at row 313 (pushing existing rows 313:358 down to 314:359)
# and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(313).Insert()

$ws.Cells.Item(313, 1).Value = 4
$ws.Cells.Item(313, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(313, 3).Value = "Los Lagos"
$ws.Cells.Item(313, 4).Value = 44984
$ws.Cells.Item(313, 5).Value = 10
$ws.Cells.Item(313, 6).Value = 100112032
$ws.Cells.Item(313, 7).Value = "Zapallo italiano"
$ws.Cells.Item(313, 8).Value = "Sin especificar"
$ws.Cells.Item(313, 9).Value = "Primera"
$ws.Cells.Item(313, 10).Value = 50
$ws.Cells.Item(313, 11).Value = 13000
$ws.Cells.Item(313, 12).Value = 13000
$ws.Cells.Item(313, 13).Value = 13000
$ws.Cells.Item(313, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(313, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(313, 16).Value = 260
$ws.Cells.Item(313, 17).Value = 50
$ws.Cells.Item(313, 18).Value = "Hortaliza"
